# Apply updates to the "yeni_degiskenler" sheet: revised OTV rate values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("yeni_degiskenler")

$ws.Range("B5").Value = 0.75
$ws.Range("B6").Value = 1.25
$ws.Range("B7").Value = 1.45
$ws.Range("B8").Value = 2.15

# Update the active selection on this sheet to match the saved view state
$ws.Activate()
$ws.Range("I10").Select()
